$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in A1:B24 with new data
$ws.Range("A1").Value = 0.00210060577194242
$ws.Range("B1").Value = 0.191068397422842
$ws.Range("A2").Value = 0.00310141725037511
$ws.Range("B2").Value = 0.0825546652729856
$ws.Range("A3").Value = 0.00147286072276185
$ws.Range("B3").Value = -0.0438553067561367
$ws.Range("A4").Value = 0.00297219028221427
$ws.Range("B4").Value = -0.070004137256296
$ws.Range("A5").Value = -0.00581863717254974
$ws.Range("B5").Value = 0.557640391800794
$ws.Range("A6").Value = 0.00146081913767116
$ws.Range("B6").Value = 0.222872867949528
$ws.Range("A7").Value = 0.00304466000818034
$ws.Range("B7").Value = 0.0872129300713072
$ws.Range("A8").Value = 0.00178448623136325
$ws.Range("B8").Value = -0.0586955531400221
$ws.Range("A9").Value = 0.00334203879907175
$ws.Range("B9").Value = -0.0854104596353763
$ws.Range("A10").Value = -0.00617119177249405
$ws.Range("B10").Value = 0.577032262183724
$ws.Range("A11").Value = 0.00227888014554518
$ws.Range("B11").Value = 0.181882213284571
$ws.Range("A12").Value = 0.00172691530189169
$ws.Range("B12").Value = 0.164076750457228
$ws.Range("A13").Value = 0.00167648269971471
$ws.Range("B13").Value = -0.051912802963853
$ws.Range("A14").Value = 0.00361339340147322
$ws.Range("B14").Value = -0.100301948901737
$ws.Range("A15").Value = -0.00572940391373669
$ws.Range("B15").Value = 0.555957925866079
$ws.Range("A16").Value = 0.00259275895280571
$ws.Range("B16").Value = 0.171033505646922
$ws.Range("A17").Value = 0.0023527967369067
$ws.Range("B17").Value = 0.121889714538855
$ws.Range("A18").Value = 0.00128508269032611
$ws.Range("B18").Value = -0.0355174157389726
$ws.Range("A19").Value = 0.00339788696290156
$ws.Range("B19").Value = -0.0920599363046228
$ws.Range("A20").Value = -0.00611960667305488
$ws.Range("B20").Value = 0.576400887889101
$ws.Range("A21").Value = 0.00162230388971332
$ws.Range("B21").Value = 0.00215043442410712
$ws.Range("A22").Value = -0.000820198240398804
$ws.Range("B22").Value = 0.0807626844727466
$ws.Range("A23").Value = 0.00458986615896255
$ws.Range("B23").Value = -0.131067075656968
$ws.Range("A24").Value = 0.0129867619495103
$ws.Range("B24").Value = -0.0835475264778341

# Remove rows 25:27 which are no longer part of the data
$ws.Range("A25:B27").Clear()
